# Re-ran "resolve" and "classify+summarise" steps after changes to the
# mapping file. This updates the SoIB_summaries workbook for West Bengal:
#   - "Range Status" sheet: species counts/percentages collapsed to 0
#     (no species resolved under the updated mapping), percentage column
#     cleared since it is no longer defined.
#   - "Species qualification" sheet: Range Analysis species count reset to 0.
#   - "High Priority break-up" sheet: the former "Range" row was dropped,
#     and the IUCN row's figures were recomputed; the old last (IUCN) row
#     is removed entirely.

$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet ---
$ws2 = $wb.Worksheets.Item("Range Status")
$ws2.Range("B2:B7").Value = 0
$ws2.Range("C2:C7").ClearContents()

# --- "Species qualification" sheet ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("B5").Value = 0

# --- "High Priority break-up" sheet ---
$ws5 = $wb.Worksheets.Item("High Priority break-up")
$ws5.Range("E2").Value = 12.1

$ws5.Range("A3").Value = "IUCN"
$ws5.Range("B3").Value = 29
$ws5.Range("C3").Value = 87.90000000000001
$ws5.Range("D3").Value = 29
$ws5.Range("E3").Value = 87.90000000000001

# The old row 4 (previous IUCN row) is no longer present; remove it so the
# sheet's dimension shrinks back to A1:E3.
$ws5.Rows("4:4").Delete()
